$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workshop 1")

# Fill in the newly-answered survey rows (columns P..S hold Q1..Q4 answers)
# Row 16
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(16, 17).Value = 0
$ws.Cells.Item(16, 18).Value = 0
$ws.Cells.Item(16, 19).Value = 0

# Row 22
$ws.Cells.Item(22, 16).Value = 1
$ws.Cells.Item(22, 17).Value = 1
$ws.Cells.Item(22, 18).Value = 0
$ws.Cells.Item(22, 19).Value = 1

# Row 28
$ws.Cells.Item(28, 16).Value = 1
$ws.Cells.Item(28, 17).Value = 1
$ws.Cells.Item(28, 18).Value = 0
$ws.Cells.Item(28, 19).Value = 0

# Update the selection / scroll position to match the edited workbook
$ws.Range("W29").Select()
